$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("B2").Value = 144298771.87442875
$ws.Range("C2").Value = 162327203.5651798
$ws.Range("D2").Value = 180355635.25593108
$ws.Range("E2").Value = 198384066.94668207
$ws.Range("F2").Value = 216412498.63743365

$ws.Range("B3").Value = 258787620.12374258
$ws.Range("C3").Value = 276816051.81449366
$ws.Range("D3").Value = 294844483.5052449
$ws.Range("E3").Value = 312872915.19599587
$ws.Range("F3").Value = 330901346.8867475

$ws.Range("B4").Value = 488005622.4446769
$ws.Range("C4").Value = 506034054.13542783
$ws.Range("D4").Value = 524062485.82617915
$ws.Range("E4").Value = 542090917.5169301
$ws.Range("F4").Value = 560119349.2076817

$ws.Range("B5").Value = 763490682.9992691
$ws.Range("C5").Value = 781519114.6900202
$ws.Range("D5").Value = 799547546.3807715
$ws.Range("E5").Value = 817575978.0715225
$ws.Range("F5").Value = 835604409.7622739
